$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1829.7097
$ws.Range("I98").Value = 2005.0834
$ws.Range("J98").Value = 1228.4286
$ws.Range("K98").Value = 2005.0834
$ws.Range("L98").Value = 1228.4286
$ws.Range("M98").Value = -507.0834
$ws.Range("N98").Value = -4224.4286
$ws.Range("H122").Value = 1829.7097
$ws.Range("I122").Value = 2005.0834
$ws.Range("J122").Value = 1228.4286
$ws.Range("K122").Value = 6015.2502
$ws.Range("L122").Value = 3685.2858
$ws.Range("M122").Value = -3565.2502
$ws.Range("N122").Value = -8585.2858
$ws.Range("H127").Value = 1621.4546
$ws.Range("I127").Value = 1449.3334
$ws.Range("J127").Value = 1828
$ws.Range("K127").Value = 4348.0002
$ws.Range("L127").Value = 5484
$ws.Range("M127").Value = 611.9997999999996
$ws.Range("N127").Value = -15404
$ws.Range("H132").Value = 289348.84
$ws.Range("I132").Value = 297835.6
$ws.Range("K132").Value = 893506.7999999999
$ws.Range("M132").Value = -890976.7999999999
$ws.Range("H137").Value = 4766273.5
$ws.Range("I137").Value = 7581498.5
$ws.Range("K137").Value = 22744495.5
$ws.Range("M137").Value = -22741945.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J4").Value = 200
$ws.Range("L4").Value = 200
$ws.Range("N4").Value = -432
$ws.Range("H11").Value = 3344002
$ws.Range("I11").Value = 3344002
$ws.Range("K11").Value = 3344002
$ws.Range("M11").Value = -3343858
$ws.Range("H74").Value = 17654810
$ws.Range("I74").Value = 30000684
$ws.Range("J74").Value = 17847.428
$ws.Range("K74").Value = 30000684
$ws.Range("L74").Value = 17847.428
$ws.Range("M74").Value = -29999810
$ws.Range("N74").Value = -19595.428
$ws.Range("H77").Value = 17654810
$ws.Range("I77").Value = 30000684
$ws.Range("J77").Value = 17847.428
$ws.Range("K77").Value = 150003420
$ws.Range("L77").Value = 89237.14
$ws.Range("M77").Value = -149999052
$ws.Range("N77").Value = -97973.14
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1685.4546
$ws.Range("I86").Value = 1537.7778
$ws.Range("J86").Value = 2350
$ws.Range("K86").Value = 1537.7778
$ws.Range("L86").Value = 2350
$ws.Range("M86").Value = -414.7778000000001
$ws.Range("N86").Value = -4596
$ws.Range("H89").Value = 1685.4546
$ws.Range("I89").Value = 1537.7778
$ws.Range("J89").Value = 2350
$ws.Range("K89").Value = 7688.889
$ws.Range("L89").Value = 11750
$ws.Range("M89").Value = -2072.889
$ws.Range("N89").Value = -22982
$ws.Range("H134").Value = 59726.21
$ws.Range("I134").Value = 85916.766
$ws.Range("J134").Value = 2980
$ws.Range("K134").Value = 257750.298
$ws.Range("L134").Value = 8940
$ws.Range("M134").Value = -255215.298
$ws.Range("N134").Value = -14010
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 2480
$ws.Range("I12").Value = 2133.3333
$ws.Range("K12").Value = 2133.3333
$ws.Range("M12").Value = -1963.3333
$ws.Range("H31").Value = 2407.75
$ws.Range("I31").Value = 1628.6666
$ws.Range("J31").Value = 3045.182
$ws.Range("K31").Value = 1628.6666
$ws.Range("L31").Value = 3045.182
$ws.Range("M31").Value = -1333.6666
$ws.Range("N31").Value = -3635.182
$ws.Range("H34").Value = 2407.75
$ws.Range("I34").Value = 1628.6666
$ws.Range("J34").Value = 3045.182
$ws.Range("K34").Value = 1628.6666
$ws.Range("L34").Value = 3045.182
$ws.Range("M34").Value = -1426.6666
$ws.Range("N34").Value = -3449.182
$ws.Range("H58").Value = 2962.95
$ws.Range("I58").Value = 1653.3334
$ws.Range("J58").Value = 3748.72
$ws.Range("K58").Value = 1653.3334
$ws.Range("L58").Value = 3748.72
$ws.Range("M58").Value = -1450.3334
$ws.Range("N58").Value = -4154.719999999999
$ws.Range("H134").Value = 1923.3158
$ws.Range("I134").Value = 1828.3334
$ws.Range("J134").Value = 2086.1428
$ws.Range("K134").Value = 5485.0002
$ws.Range("L134").Value = 6258.428400000001
$ws.Range("M134").Value = -2950.0002
$ws.Range("N134").Value = -11328.4284
$ws.Range("H136").Value = 2962.95
$ws.Range("I136").Value = 1653.3334
$ws.Range("J136").Value = 3748.72
$ws.Range("K136").Value = 4960.0002
$ws.Range("L136").Value = 11246.16
$ws.Range("M136").Value = -2410.0002
$ws.Range("N136").Value = -16346.16
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 952
$ws.Range("I3").Value = 952
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2856
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -2744
$ws.Range("N3").ClearContents()
$ws.Range("H5").Value = 522.2941
$ws.Range("I5").Value = 499.91666
$ws.Range("J5").Value = 576
$ws.Range("K5").Value = 1499.74998
$ws.Range("L5").Value = 1728
$ws.Range("M5").Value = -1387.74998
$ws.Range("N5").Value = -1952
$ws.Range("H6").Value = 70
$ws.Range("I6").Value = 70
$ws.Range("K6").Value = 210
$ws.Range("M6").Value = -97
$ws.Range("H86").Value = 1333.3334
$ws.Range("I86").Value = 1333.3334
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 4000.0002
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -2814.0002
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 1333.3334
$ws.Range("I89").Value = 1333.3334
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 12000.0006
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -6072.000599999999
$ws.Range("N89").ClearContents()
$ws.Range("H98").Value = 691.8421
$ws.Range("I98").Value = 546.5833
$ws.Range("J98").Value = 940.8570999999999
$ws.Range("K98").Value = 1639.7499
$ws.Range("L98").Value = 2822.5713
$ws.Range("M98").Value = -141.7499
$ws.Range("N98").Value = -5818.5713
$ws.Range("H131").Value = 1236263.5
$ws.Range("J131").Value = 1353061.8
$ws.Range("L131").Value = 4059185.4
$ws.Range("N131").Value = -4069265.4
$ws.Range("H135").Value = 522.2941
$ws.Range("I135").Value = 499.91666
$ws.Range("J135").Value = 576
$ws.Range("K135").Value = 4499.24994
$ws.Range("L135").Value = 5184
$ws.Range("M135").Value = -1964.24994
$ws.Range("N135").Value = -10254
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 7008
$ws.Range("I17").Value = 7008
$ws.Range("K17").Value = 7008
$ws.Range("M17").Value = -6840
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H70").Value = 4768.84
$ws.Range("I70").Value = 4894.7856
$ws.Range("J70").Value = 4608.5454
$ws.Range("K70").Value = 4894.7856
$ws.Range("L70").Value = 4608.5454
$ws.Range("M70").Value = -4624.7856
$ws.Range("N70").Value = -5148.5454
$ws.Range("H73").Value = 4768.84
$ws.Range("I73").Value = 4894.7856
$ws.Range("J73").Value = 4608.5454
$ws.Range("K73").Value = 4894.7856
$ws.Range("L73").Value = 4608.5454
$ws.Range("M73").Value = -3958.7856
$ws.Range("N73").Value = -6480.5454
$ws.Range("H122").Value = 5435.857
$ws.Range("I122").Value = 7210.25
$ws.Range("J122").Value = 3070
$ws.Range("K122").Value = 21630.75
$ws.Range("L122").Value = 9210
$ws.Range("M122").Value = -19180.75
$ws.Range("N122").Value = -14110
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 501.53845
$ws.Range("I22").Value = 598.5714
$ws.Range("J22").Value = 388.33334
$ws.Range("K22").Value = 598.5714
$ws.Range("L22").Value = 388.33334
$ws.Range("M22").Value = -303.5714
$ws.Range("N22").Value = -978.33334
$ws.Range("H27").Value = 501.53845
$ws.Range("I27").Value = 598.5714
$ws.Range("J27").Value = 388.33334
$ws.Range("K27").Value = 598.5714
$ws.Range("L27").Value = 388.33334
$ws.Range("M27").Value = -491.5714
$ws.Range("N27").Value = -602.33334
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 50000
$ws.Range("J46").Value = 50000
$ws.Range("L46").Value = 50000
$ws.Range("N46").Value = -50462
$ws.Range("H126").Value = 837.4286
$ws.Range("I126").Value = 863.9524
$ws.Range("K126").Value = 2591.8572
$ws.Range("M126").Value = -121.8571999999999
$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 150000
$ws.Range("N134").Value = -155070
$ws.Range("H136").Value = 32664.37
$ws.Range("I136").Value = 9001.833000000001
$ws.Range("J136").Value = 73228.71000000001
$ws.Range("K136").Value = 27005.499
$ws.Range("L136").Value = 219686.13
$ws.Range("M136").Value = -24455.499
$ws.Range("N136").Value = -224786.13

Write-Host "Applied 222 cell changes across 8 sheets"